$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.275.43"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.128.79"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.52"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.56"
$ws.Range("E6").Value = "  -1.35%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.123.35"
$ws.Range("E8").Value = "  +1.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("E10").Value = "  +1.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("E11").Value = "  +3.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  +3.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.07"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.653.13"
$ws.Range("E15").Value = "  +1.56%  "

$ws.Range("E16").Value = "  +3.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.230.71"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.116.71"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.04"
$ws.Range("E20").Value = "  +1.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.709"
$ws.Range("E22").Value = "  +1.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.65"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("E24").Value = "  +2.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.35"
$ws.Range("E25").Value = "  -0.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.34"
$ws.Range("E28").Value = "  +0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  +8.36%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.114"
$ws.Range("E30").Value = "  +2.34%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  -4.38%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.82"
$ws.Range("E33").Value = "  +3.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("E34").Value = "  -2.89%  "

$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0767"
$ws.Range("E36").Value = "  +5.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.94"
$ws.Range("E37").Value = "  +0.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.28"
$ws.Range("E38").Value = "  +0.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  +4.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "443.93"
$ws.Range("E40").Value = "  -2.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0392"
$ws.Range("E41").Value = "  +0.48%  "

$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.18"
$ws.Range("E43").Value = "  -1.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.849.91"
$ws.Range("E44").Value = "  +1.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.259"
$ws.Range("E45").Value = "  -1.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.91"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.75"
$ws.Range("E51").Value = "  +1.75%  "
